# Swap the colour scheme that is bound to the slide master (ppt/theme/theme1.xml,
# currently the "Integral" / "Red Violet" theme) so that it becomes the standard
# "Office Theme" colour scheme (the one that, before this edit, only lived in
# ppt/theme/theme2.xml, used by the notes master).
$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # dk1
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1
$tcs.Item(3).RGB  = 0x6A5444   # dk2
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2
$tcs.Item(5).RGB  = 0xD59B5B   # accent1
$tcs.Item(6).RGB  = 0x317DED   # accent2
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3
$tcs.Item(8).RGB  = 0x00C0FF   # accent4
$tcs.Item(9).RGB  = 0xC47244   # accent5
$tcs.Item(10).RGB = 0x47AD70   # accent6
$tcs.Item(11).RGB = 0xC16305   # hlink
$tcs.Item(12).RGB = 0x724F95   # folHlink
